$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix/expand god-name cells with their ability descriptions
$ws.Range("B2").Value = "Posiedon: Opponents God ability cannot be used"
$ws.Range("D2").Value = "Neptune: put opponents creatures to sleep that matches your creatures"
$ws.Range("D4").Value = 'Venus: One card goes back to your hand if it "dies"'
$ws.Range("B5").Value = "Zeus: draw a card from opponents hand and play it"
$ws.Range("D5").Value = "Jupiter: Opponents hand gets shuffled and they can't see any of their cards"
$ws.Range("D6").Value = "Pluto: One card from opponents discard goes to your hand"

# New entry added further down the sheet
$ws.Range("B11").Value = "Copy opponents card"
